$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A13,A19")
foreach ($a in $r.Areas) { $a.Font.ThemeColor = 2 }
foreach ($a in $r.Areas) { $a.Interior.ThemeColor = 1 }
